$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# New row: Worksheet | autofilter | excel-worksheet-auto-filter | addAutoFilter
$row1 = $lo.ListRows.Add()
$row1.Range.Item(1,1).Value = "Worksheet"
$row1.Range.Item(1,2).Value = "autofilter"
$row1.Range.Item(1,4).Value = "addAutoFilter"
$row1.Range.Item(1,3).Value = "excel-worksheet-auto-filter"

# New row: AutoFilter | apply | excel-worksheet-auto-filter | addAutoFilter
$row2 = $lo.ListRows.Add()
$row2.Range.Item(1,1).Value = "AutoFilter"
$row2.Range.Item(1,2).Value = "apply"
$row2.Range.Item(1,4).Value = "addAutoFilter"
$row2.Range.Item(1,3).Value = "excel-worksheet-auto-filter"

# Match the final saved selection from the authored edit
$ws.Range("O178").Select()
